$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert 4 new rows before row 21 (pushes existing Day5..Day13 content down to rows 25-48)
$ws.Rows("21:24").Insert()

# Fill in the new "Day 5" block (Revision HTML and CSS)
$ws.Range("A21").Value = "Day 5"
$ws.Range("B21").Value = "Revision HTML and CSS"
$ws.Range("C21").Value = "HTML: Introduction, Sementic tags, Emmet Abbreviations"
$ws.Range("C22").Value = "CSS: Box Model, Position, Float, Flex"
$ws.Range("C23").Value = "CSS: Grid, Transition, Animation"
$ws.Range("C24").Value = "CSS: Media Query"

# Merge the Day / Module Name columns for the new block
$ws.Range("A21:A24").Merge()
$ws.Range("B21:B24").Merge()

# Extend the bold formatting that previously stopped at row 14 down through
# the Day 4 block and the new Day 5 (HTML/CSS) block (rows 15-24)
$ws.Range("A15:C24").Font.Bold = $true

# Restore selection / scroll position similar to where the author left off editing
$ws.Application.GoTo($ws.Range("B25"), $true)
$ws.Range("B25:B28").Select()
$ws.Application.ActiveWindow.ScrollRow = 13

Write-Host "Edit applied successfully"
